$d = $word.ActiveDocument
$d.Content.Find.Execute("2019/01/16 10:01:54 - Lost user content zone1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2019/02/08 16:15:54 - Lost user content zone1", 2)
